$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Sheet "VENTAS POR GRUPO": D5 sale value, D23 running count label
$wsVentasPorGrupo.Range("D5").Value = 915.84
$wsVentasPorGrupo.Range("D23").Value = "1 de 21"

# Sheet "VENTA MENSUAL": F5 sale value for agosto, F23 column total
$wsVentaMensual.Range("F5").Value = 915.84
$wsVentaMensual.Range("F23").Value = 7471.97

# Sheet "CUMPLIMIENTO MENSUAL": row 3 (240X80 PORCELANATO) and row 19 (TOTAL)
$wsCumplimientoMensual.Range("D3").Value = 915.84
$wsCumplimientoMensual.Range("E3").Value = 3252.23156573679
$wsCumplimientoMensual.Range("F3").Value = 0.2197275132050443

$wsCumplimientoMensual.Range("D19").Value = 7471.97
$wsCumplimientoMensual.Range("E19").Value = 51916.25762291768
$wsCumplimientoMensual.Range("F19").Value = 0.1258156759188516
